$wb = $excel.ActiveWorkbook

# Update version string on the isa_template sheet (B4: 1.0.0 -> 1.0.1)
$tmpl = $wb.Worksheets.Item("isa_template")
$tmpl.Range("B4").Value = "1.0.1"

# Fill in example values in row 2 of the "MS" sheet
$ms = $wb.Worksheets.Item("MS")
$ms.Range("B2").Value = "Mass Spectrometry"
$ms.Range("C2").Value = "NCIT"
$ms.Range("D2").Value = "http://purl.obolibrary.org/obo/NCIT_C17156"
$ms.Range("F2").Value = "positive scan"
$ms.Range("G2").Value = "MS"
$ms.Range("H2").Value = "http://purl.obolibrary.org/obo/MS_1000130"
$ms.Range("I2").Value = "100-1000"
$ms.Range("L2").Value = "Bruker micrOTOF-Q II"
$ms.Range("O2").Value = "electrospray ionization"
$ms.Range("P2").Value = "MS"
$ms.Range("Q2").Value = "http://purl.obolibrary.org/obo/MS_1000073"
$ms.Range("R2").Value = "triple quadrupole"
